# FrameDesign.xlsx — "Change frame design + Trasmitter CodeV1"
#
# A new "src 8 bits" column is inserted into the frame-design header row
# between the existing "dst 8 bits" (B) and "frameNo 1 bit" (C) columns.
# The previous C:F columns (frameNo 1 bit / data 10 bits / parity 1 bit /
# endFlag 1 bit, plus the "1 for more" / "0 for last" notes in rows 2-3)
# all shift one column to the right (C->D, D->E, E->F, F->G), the brand
# new column C is filled in with "src 8 bits", and the new trailing
# column G is given its own width. The active selection ends up on C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the header row (and the two note cells in F2/F3) one column to the
# right. Cut (not Copy) so the vacated source cell's content/shared-string
# moves rather than being duplicated, working right-to-left so we never
# clobber a column before it has been moved.
$ws.Range("F1").Cut($ws.Range("G1"))
$ws.Range("E1").Cut($ws.Range("F1"))
$ws.Range("D1").Cut($ws.Range("E1"))
$ws.Range("C1").Cut($ws.Range("D1"))

$ws.Range("F2").Cut($ws.Range("G2"))
$ws.Range("F3").Cut($ws.Range("G3"))

# New "src 8 bits" header, dropped into the now-empty (but still styled)
# C1 cell vacated by the cut above.
$ws.Range("C1").Value = "src 8 bits"

# New column G gets its own explicit width (matches the other header
# columns' custom widths).
$ws.Columns.Item(7).ColumnWidth = 12

# Final selection, matching the saved sheet view.
$ws.Range("C5").Select()
